# adding averages and more checks
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Training Dashboard (sheet 1): refresh "PERIOD TO EXPIRE" (H) and
#    "LAST UPDATE" (I) for every data row (3-23). The "LAST UPDATE" date
#    moves from 08-Sep-2025 to 16-Sep-2025, and "PERIOD TO EXPIRE" drops by
#    the same 8 days.
# ---------------------------------------------------------------------------
$wsTraining = $wb.Worksheets.Item("Training Dashboard")

$periodToExpire = @{
    3  = 402
    4  = 323
    5  = 325
    6  = 363
    7  = 352
    8  = 688
    9  = 406
    10 = 385
    11 = 365
    12 = 402
    13 = 364
    14 = 388
    15 = 392
    16 = 406
    17 = 405
    18 = 362
    19 = 85
    20 = 174
    21 = 177
    22 = 189
    23 = 232
}

foreach ($row in $periodToExpire.Keys) {
    $wsTraining.Cells.Item($row, 8).Value = $periodToExpire[$row]
    # Leading apostrophe keeps this a literal text value (matches the
    # existing inline-string "LAST UPDATE" cells) instead of Excel's
    # automatic date-serial conversion.
    $wsTraining.Cells.Item($row, 9).Value = "'16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 2) Exam Dashboard (sheet 2): widen the COMMENTS column and update the
#    first remark from "OK" to "date is valid".
# ---------------------------------------------------------------------------
$wsExam = $wb.Worksheets.Item("Exam Dashboard")

$wsExam.Columns.Item(5).ColumnWidth = 14.1667
$wsExam.Range("E3").Value = "date is valid"

# ---------------------------------------------------------------------------
# 3) Recolor the dashboard titles and header rows to bold white text
#    (reusing one shared bold font instead of the old two-font setup).
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $usedCols = $ws.UsedRange.Columns.Count

    $title = $ws.Range("A1")
    $title.Font.Size = 11
    $title.Font.Bold = $true
    $title.Font.Color = 16777215

    $header = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(2, $usedCols))
    $header.Font.Bold = $true
    $header.Font.Color = 16777215
}
